# Add two new worksheets (Sheet4, Sheet5) at the end of the workbook,
# populate them with QR-code label data, and make Sheet4 the active tab.

$wb = $excel.ActiveWorkbook

# --- Sheet4 --------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$s4 = $wb.Worksheets.Add($null, $lastSheet)

$s4.Range("A1").Value = "QRNUM"
$s4.Range("A2").Value = "DGR OSA TO1"
$s4.Range("A3").Value = "DGR OSA TO2"
$s4.Range("A4").Value = "DGR OSA TO3"
$s4.Range("A5").Value = "DGR OSA TO4"
$s4.Range("A6").Value = "DGR OSA TO5"
$s4.Range("A7").Value = "DGR OSA TO6"
$s4.Range("A8").Value = "DGR OSA TO7"
$s4.Range("A9").Value = "DGR OSA TO8"
$s4.Range("A10").Value = "DGR OSA TO9"
$s4.Range("A11").Value = "DBT OSA G09"
$s4.Range("A12").Value = "DBT OSA I07"
$s4.Range("A13").Value = "DBT OSA I08"

$s4.Columns.Item(1).ColumnWidth = 16.666666666666668

# --- Sheet5 ----------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$s5 = $wb.Worksheets.Add($null, $lastSheet2)

$s5.Range("A1").Value = "QRNUM"
$s5.Range("A2").Value = "DRB OPR B02"
$s5.Range("A3").Value = "DRB OPR B03"

$s5.Columns.Item(1).ColumnWidth = 12

# --- selections / active sheet ---------------------------------------------
$s5.Range("A2:A3").Select()
$s4.Range("A2:A13").Select()
$s4.Activate()
